# AFDP-2522 ensure case files rules work before queue has been assigned
#
# The "queue" variable on a case file may be null before the case has been
# routed to a queue at all. The rule conditions used a plain property
# access (queue.name) which blows up with a NullPointerException when
# queue is null. Switch to the null-safe navigation operator (queue?.name)
# so the rules keep working even when no queue has been assigned yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set Billing Enter Date / Nullify Billing Enter Date / Set Hold Enter Date /
# Nullify Hold Enter Date condition expressions (column C, rows 28-31).
$ws.Range("C28").Value = "queue?.name == 'Billing' && billingEnterDate == null"
$ws.Range("C29").Value = "queue?.name != 'Billing'"
$ws.Range("C30").Value = "queue?.name == 'Hold' && holdEnterDate == null"
$ws.Range("C31").Value = "queue?.name != 'Hold'"

# Tidy up the trailing blank rows: the leftover styled-but-empty cells in
# row 32 (B32:D32) and the now-redundant blank row 33 are removed.
$ws.Range("B32:D32").Clear()
$ws.Rows.Item(33).Delete()

# Leave the cursor where the author left it after making the edit.
$ws.Range("B26").Select()
